$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "door_nodes": the row that held the stale 31030/2112 pair (the old
# row 6, id=5) is being removed; its real data (24373/24369/...) is copied
# into the two rows that used to carry the placeholder values (rows 4 & 5,
# ids 3 & 4). Deleting that row then shifts the two rows below it up by
# one (they already carried the correct ids 6 & 7, so nothing else to fix).
# -----------------------------------------------------------------------
$doors = $wb.Worksheets.Item("door_nodes")

foreach ($r in 4, 5) {
    $doors.Range("C$r").Value = 24373.0
    $doors.Range("D$r").Value = 24369.0
    $doors.Range("E$r").Value = 2815.0
    $doors.Range("F$r").Value = 3182.0
    $doors.Range("G$r").Value = 0.0
    $doors.Range("H$r").Value = 2816.0
    $doors.Range("I$r").Value = 3182.0
    $doors.Range("J$r").Value = 0.0
    $doors.Range("K$r").Value = 2815.0
    $doors.Range("L$r").Value = 3182.0
    $doors.Range("M$r").Value = 0.0
    $doors.Range("N$r").Value = 2816.0
    $doors.Range("O$r").Value = 3182.0
}

$doors.Rows.Item(6).Delete()

# -----------------------------------------------------------------------
# Sheet "object_nodes": the importer can now match an object purely by
# name instead of by object_id, so row 2 (the "Climb-down" entry) becomes
# a name-matched "Ladder" row, and two new id-matched rows (copies of the
# door's open/close tile data, id 2112) are appended for that same ladder.
# -----------------------------------------------------------------------
$objects = $wb.Worksheets.Item("object_nodes")

$objects.Range("B2").Value = "name"
$objects.Range("C2").ClearContents()
$objects.Range("C2").Copy()
$objects.Range("D2").PasteSpecial(-4122)
$objects.Range("D2").Value = "Ladder"

foreach ($col in "A","B","C","E","F","G","H","I","J","K","L","M","N","O","P","T") {
    $objects.Range("$col`2").Copy()
    $objects.Range("$col`4").PasteSpecial(-4122)
    $objects.Range("$col`5").PasteSpecial(-4122)
}

$objects.Range("A4").Value = 3.0
$objects.Range("B4").Value = "id"
$objects.Range("C4").Value = 2112.0
$objects.Range("E4").Value = "Open"
$objects.Range("F4").Value = 3045.0
$objects.Range("G4").Value = 3047.0
$objects.Range("H4").Value = 9755.0
$objects.Range("I4").Value = 9756.0
$objects.Range("J4").Value = 0.0
$objects.Range("K4").Value = 3045.0
$objects.Range("L4").Value = 3047.0
$objects.Range("M4").Value = 9757.0
$objects.Range("N4").Value = 9758.0
$objects.Range("O4").Value = 0.0
$objects.Range("P4").Value = 20.0
$objects.Range("T4").Value = 30.0

$objects.Range("A5").Value = 4.0
$objects.Range("B5").Value = "id"
$objects.Range("C5").Value = 2112.0
$objects.Range("E5").Value = "Open"
$objects.Range("F5").Value = 3045.0
$objects.Range("G5").Value = 3047.0
$objects.Range("H5").Value = 9757.0
$objects.Range("I5").Value = 9758.0
$objects.Range("J5").Value = 0.0
$objects.Range("K5").Value = 3045.0
$objects.Range("L5").Value = 3047.0
$objects.Range("M5").Value = 9755.0
$objects.Range("N5").Value = 9756.0
$objects.Range("O5").Value = 0.0
$objects.Range("P5").Value = 20.0
$objects.Range("T5").Value = 30.0
